$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Previous account-statement rows (16-20) held:
#   row16: PE / 90086118011978 / SEIR SOTILLO MONSERRAT / 2204 / 12800  / 1200000
#   row17: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2308 / 80000  / 2000000
#   row18: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2309 / 80000  / 2000000
#   row19: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2310 / 80000  / 2000000
#   row20: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2311 / 80000  / 2000000
#
# New data: remove the old "SEIR SOTILLO MONSERRAT" period, push it to the bottom,
# and list the "JESUS ALEJANDRO LEON MATA" periods in reverse (2311..2308) starting
# at the top:
#   row16: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2311 / 80000  / 2000000
#   row17: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2310 / 80000  / 2000000
#   row18: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2309 / 80000  / 2000000
#   row19: PPT / 1093226 / JESUS ALEJANDRO LEON MATA    / 2308 / 80000  / 2000000
#   row20: PE  / 90086118011978 / SEIR SOTILLO MONSERRAT / 2204 / 12800 / 1200000

$ws.Range("B16").Value = "PPT"
$ws.Range("C16").Value = "1093226"
$ws.Range("D16").Value = "JESUS ALEJANDRO LEON MATA"
$ws.Range("E16").Value = "2311"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

$ws.Range("B17").Value = "PPT"
$ws.Range("C17").Value = "1093226"
$ws.Range("D17").Value = "JESUS ALEJANDRO LEON MATA"
$ws.Range("E17").Value = "2310"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

$ws.Range("B18").Value = "PPT"
$ws.Range("C18").Value = "1093226"
$ws.Range("D18").Value = "JESUS ALEJANDRO LEON MATA"
$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 80000
$ws.Range("G18").Value = 2000000

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "1093226"
$ws.Range("D19").Value = "JESUS ALEJANDRO LEON MATA"
$ws.Range("E19").Value = "2308"
$ws.Range("F19").Value = 80000
$ws.Range("G19").Value = 2000000

$ws.Range("B20").Value = "PE"
$ws.Range("C20").Value = "90086118011978"
$ws.Range("D20").Value = "SEIR SOTILLO MONSERRAT"
$ws.Range("E20").Value = "2204"
$ws.Range("F20").Value = 12800
$ws.Range("G20").Value = 1200000
